$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 52.049028
$ws.Range("H2").Value = 156.147084
$ws.Range("I2").Value = 0.7208330343078339
$ws.Range("J2").Value = 0.7208330343078339
$ws.Range("M2").Value = 36.923013
$ws.Range("N2").Value = 110.769039
$ws.Range("O2").Value = 0.7437819354528793
$ws.Range("P2").Value = 0.7437819354528794
$ws.Range("Q2").Value = 1921.806937481364
$ws.Range("R2").Value = 17296.26243733228
$ws.Range("S2").Value = 0.5361425893958524
$ws.Range("T2").Value = 0.5361425893958525

# Row 3
$ws.Range("G3").Value = 52.049028
$ws.Range("H3").Value = 156.147084
$ws.Range("I3").Value = 0.7208330343078339
$ws.Range("J3").Value = 0.7208330343078339
$ws.Range("O3").Value = 0.1364233939221953
$ws.Range("P3").Value = 0.1364233939221953
$ws.Range("Q3").Value = 352.49501551122
$ws.Range("R3").Value = 3172.45513960098
$ws.Range("S3").Value = 0.09833848899150896
$ws.Range("T3").Value = 0.09833848899150896

# Row 4
$ws.Range("G4").Value = 52.049028
$ws.Range("H4").Value = 156.147084
$ws.Range("I4").Value = 0.7208330343078339
$ws.Range("J4").Value = 0.7208330343078339
$ws.Range("O4").Value = 0.1197946706249253
$ws.Range("P4").Value = 0.1197946706249254
$ws.Range("Q4").Value = 309.529202184908
$ws.Range("R4").Value = 2785.762819664172
$ws.Range("S4").Value = 0.08635195592047247
$ws.Range("T4").Value = 0.08635195592047247

# Row 5
$ws.Range("I5").Value = 0.09317473454775864
$ws.Range("J5").Value = 0.09317473454775864
$ws.Range("M5").Value = 36.923013
$ws.Range("N5").Value = 110.769039
$ws.Range("O5").Value = 0.7437819354528793
$ws.Range("P5").Value = 0.7437819354528794
$ws.Range("Q5").Value = 248.412382243011
$ws.Range("R5").Value = 2235.711440187099
$ws.Range("S5").Value = 0.06930168439724017
$ws.Range("T5").Value = 0.06930168439724019

# Row 6
$ws.Range("I6").Value = 0.09317473454775864
$ws.Range("J6").Value = 0.09317473454775864
$ws.Range("O6").Value = 0.1364233939221953
$ws.Range("P6").Value = 0.1364233939221953
$ws.Range("S6").Value = 0.01271121351480486
$ws.Range("T6").Value = 0.01271121351480486

# Row 7
$ws.Range("I7").Value = 0.09317473454775864
$ws.Range("J7").Value = 0.09317473454775864
$ws.Range("O7").Value = 0.1197946706249253
$ws.Range("P7").Value = 0.1197946706249254
$ws.Range("S7").Value = 0.0111618366357136
$ws.Range("T7").Value = 0.0111618366357136

# Row 8
$ws.Range("I8").Value = 0.1859922311444076
$ws.Range("J8").Value = 0.1859922311444076
$ws.Range("M8").Value = 36.923013
$ws.Range("N8").Value = 110.769039
$ws.Range("O8").Value = 0.7437819354528793
$ws.Range("P8").Value = 0.7437819354528794
$ws.Range("Q8").Value = 495.872335365687
$ws.Range("R8").Value = 4462.851018291183
$ws.Range("S8").Value = 0.1383376616597868
$ws.Range("T8").Value = 0.1383376616597868

# Row 9
$ws.Range("I9").Value = 0.1859922311444076
$ws.Range("J9").Value = 0.1859922311444076
$ws.Range("O9").Value = 0.1364233939221953
$ws.Range("P9").Value = 0.1364233939221953
$ws.Range("Q9").Value = 90.952177941135
$ws.Range("R9").Value = 818.569601470215
$ws.Range("S9").Value = 0.02537369141588152
$ws.Range("T9").Value = 0.02537369141588152

# Row 10
$ws.Range("I10").Value = 0.1859922311444076
$ws.Range("J10").Value = 0.1859922311444076
$ws.Range("O10").Value = 0.1197946706249253
$ws.Range("P10").Value = 0.1197946706249254
$ws.Range("Q10").Value = 79.865966428689
$ws.Range("R10").Value = 718.7936978582011
$ws.Range("S10").Value = 0.02228087806873929
$ws.Range("T10").Value = 0.02228087806873929
